# Changed to use different initial password; and to use a fixed final password
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Global")

# A2 held the "initial password" (j.kaplan@microfocus.com) -> use a different one
$ws.Range("A2").Value = "ppm_octane@microfocus.com"

# A5 held the "final password" (ppm_octane@microfocus.com) -> now fixed to j.kaplan@microfocus.com
$ws.Range("A5").Value = "j.kaplan@microfocus.com"

# Update the saved selection on the Global sheet to A6, without disturbing
# which sheet/tab is actually active in the workbook.
$previousActiveSheet = $wb.ActiveSheet.Name
$ws.Activate()
$ws.Range("A6").Select()
$wb.Worksheets.Item($previousActiveSheet).Activate()
